$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - value + style matching the other header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Timestamps for F2:F27, one per data row
$timestamps = @(
    "2021-10-05 13:39:38.586778",
    "2021-10-05 13:39:38.586791",
    "2021-10-05 13:39:38.586795",
    "2021-10-05 13:39:38.586798",
    "2021-10-05 13:39:38.586802",
    "2021-10-05 13:39:38.586805",
    "2021-10-05 13:39:38.586808",
    "2021-10-05 13:39:38.586811",
    "2021-10-05 13:39:38.586815",
    "2021-10-05 13:39:38.586818",
    "2021-10-05 13:39:38.586821",
    "2021-10-05 13:39:38.586825",
    "2021-10-05 13:39:38.586828",
    "2021-10-05 13:39:38.586831",
    "2021-10-05 13:39:38.586834",
    "2021-10-05 13:39:38.586837",
    "2021-10-05 13:39:38.586841",
    "2021-10-05 13:39:38.586844",
    "2021-10-05 13:39:38.586847",
    "2021-10-05 13:39:38.586850",
    "2021-10-05 13:39:38.586853",
    "2021-10-05 13:39:38.586857",
    "2021-10-05 13:39:38.586860",
    "2021-10-05 13:39:38.586863",
    "2021-10-05 13:39:38.586866",
    "2021-10-05 13:39:38.586870"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
